# Remove the paragraph describing who manages the social-media pages
# ("As redes Sociais, site são gerenciadas por Reinaldo e Clayton
# colaboradores da empresa"), which sat directly under the
# "Referências" heading. The heading paragraph and everything after
# the removed paragraph stay untouched.

$d = $word.ActiveDocument

# Locate the paragraph by a snippet of its (unique) text so this does
# not depend on a hard-coded paragraph index.
$rng = $d.Content
$found = $rng.Find.Execute("As redes So", $false, $false, $false, $false, `
                            $false, $true, 1, $false, "", 0)

if ($found) {
    # Find which paragraph in the document contains the match, then
    # grab that Paragraph object so we can remove it (text + its own
    # paragraph mark) as a single unit.
    $count = $d.Paragraphs.Count
    $target = $null
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($rng.Start -ge $p.Range.Start -and $rng.Start -lt $p.Range.End) {
            $target = $p
            break
        }
    }

    if ($target -ne $null) {
        $target.Range.Delete()
    }
}
